$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    # Force the value to be stored as text (matches workbook convention where
    # every data cell - including numeric-looking ones - is inline string text).
    if ($text -match '^-?[0-9]') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
    $range.Style = "Normal"
}

# --- USERS sheet ---
$usersWs = $wb.Worksheets.Item("USERS")
Set-TextValue $usersWs.Range("B2") "hello"
Set-TextValue $usersWs.Range("C2") "hello"
Set-TextValue $usersWs.Range("B4") "3"
Set-TextValue $usersWs.Range("C4") "3"
Set-TextValue $usersWs.Range("B5") "q"
Set-TextValue $usersWs.Range("C5") "q"

# --- PARAMETERS sheet ---
$paramsWs = $wb.Worksheets.Item("PARAMETERS")

# Row 2: AOO mode, tuned-down parameters
Set-TextValue $paramsWs.Range("B2") "AOO"
Set-TextValue $paramsWs.Range("C2") "50"
Set-TextValue $paramsWs.Range("D2") "50"
Set-TextValue $paramsWs.Range("F2") "off"
Set-TextValue $paramsWs.Range("H2") "0.05"
Set-TextValue $paramsWs.Range("J2") "500"
Set-TextValue $paramsWs.Range("K2") "500"
Set-TextValue $paramsWs.Range("L2") "V-Low"
Set-TextValue $paramsWs.Range("M2") "50"
Set-TextValue $paramsWs.Range("N2") "5"

# Row 3: AAI mode
Set-TextValue $paramsWs.Range("B3") "AAI"
Set-TextValue $paramsWs.Range("F3") "3.5"

# Row 4: add Mode column value
Set-TextValue $paramsWs.Range("B4") "AOO"

# Row 5: new full row of parameters (AOO mode defaults)
Set-TextValue $paramsWs.Range("B5") "AOO"
Set-TextValue $paramsWs.Range("C5") "60"
Set-TextValue $paramsWs.Range("D5") "120"
Set-TextValue $paramsWs.Range("E5") "150"
Set-TextValue $paramsWs.Range("F5") "3.5"
Set-TextValue $paramsWs.Range("G5") "3.5"
Set-TextValue $paramsWs.Range("H5") "0.4"
Set-TextValue $paramsWs.Range("I5") "0.5"
Set-TextValue $paramsWs.Range("J5") "320"
Set-TextValue $paramsWs.Range("K5") "250"
Set-TextValue $paramsWs.Range("L5") "Med"
Set-TextValue $paramsWs.Range("M5") "30"
Set-TextValue $paramsWs.Range("N5") "8"
Set-TextValue $paramsWs.Range("O5") "5"

$wb.Save()
